# Upload til SharePoint Testsite
# Insert three new case/email rows into the "Indsender emails" sheet while
# keeping the existing rows (and their formatting) intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row at the very top of the data (becomes row 2; everything below
# shifts down by one). ---
$ws.Rows.Item(2).Insert()
# Excel's default row Insert copies the format of the row ABOVE the
# insertion point. For row 2 that is the bold/blue header row, so reset the
# new row back to the plain bordered style used by the rest of the data
# rows before stamping the border back on.
$ws.Range("A2:B2").Style = "Normal"
$ws.Range("A2:B2").Borders.LineStyle = 1
$ws.Range("A2").Value = "S2026-3001"
$ws.Range("B2").Value = "ufda@aarhus.dk"

# --- New row inserted right after "S2025-82487" / before "S2025-82157".
# After the first insertion above, that position is row 14. ---
$ws.Rows.Item(14).Insert()
$ws.Range("A14:B14").Style = "Normal"
$ws.Range("A14:B14").Borders.LineStyle = 1
$ws.Range("A14").Value = "S2025-82367"
$ws.Range("B14").Value = "uffe@fredens.net"

# --- New row inserted right after "S2025-71626" / before "S2025-63378".
# After the two prior insertions, that position is row 32. ---
$ws.Rows.Item(32).Insert()
$ws.Range("A32:B32").Style = "Normal"
$ws.Range("A32:B32").Borders.LineStyle = 1
$ws.Range("A32").Value = "S2025-68291"
$ws.Range("B32").Value = "line.bak.elleskov@danbolig.dk"
